$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 World Cup")
$ws.Activate()

# Fill in the final group-stage match scores (row 51-54)
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 1

$ws.Range("F52").Value = 1
$ws.Range("G52").Value = 2

$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 1

$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 1

# Update the visible selection/scroll position to match the new view state
$ws.Range("F49").Select()
